$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback datetimes on row 5
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-02-26 07:08:18"
$wsZh.Range("G5").Value = "2016-02-26 07:09:05"

# de-de sheet: update Correspond Handoff/Handback datetimes on row 5
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-02-26 07:08:33"
$wsDe.Range("G5").Value = "2016-02-26 07:09:28"
